# Update the JSON payload strings in column B (rows 2-42) to add the
# 'userType': 'TENANT' key, matching the commit "Existing API updated with
# latest payload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value2
    if ($current -ne $null -and $current -match "^\{'email': '([^']*)'\}$") {
        $email = $matches[1]
        $cell.Value = "{'email': '$email', 'userType': 'TENANT'}"
    }
}
